$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = 0.1124019449155688
$ws.Range("C2").Value = 0.6821200353019861
$ws.Range("D2").Value = 1.209983865900432
$ws.Range("E2").Value = 1.099992666293931
$ws.Range("F2").Value = 1.105344008082418
$ws.Range("G2").Value = 50

$ws.Range("B3").Value = 0.1960116252932615
$ws.Range("C3").Value = 0.6752669640606744
$ws.Range("D3").Value = 1.511496450551623
$ws.Range("E3").Value = 1.22942931905483
$ws.Range("F3").Value = 1.226547159135974
$ws.Range("G3").Value = 48

$ws.Range("B4").Value = 0.001728000735583749
$ws.Range("C4").Value = 0.6309824882362022
$ws.Range("D4").Value = 1.099648267891479
$ws.Range("E4").Value = 1.048641153060225
$ws.Range("F4").Value = 1.059736755062583
$ws.Range("G4").Value = 48

$ws.Range("B5").Value = 0.0937284794034382
$ws.Range("C5").Value = 0.6894373285789633
$ws.Range("D5").Value = 1.642876430466323
$ws.Range("E5").Value = 1.281747412896286
$ws.Range("F5").Value = 1.292135878597467
$ws.Range("G5").Value = 47

$ws.Range("B6").Value = 0.04836703510120018
$ws.Range("C6").Value = 0.6561079152665072
$ws.Range("D6").Value = 1.680400965917633
$ws.Range("E6").Value = 1.296302806414316
$ws.Range("F6").Value = 1.309714416778838
$ws.Range("G6").Value = 46

$ws.Range("B7").Value = 0.0481232454536701
$ws.Range("C7").Value = 0.6637583974150937
$ws.Range("D7").Value = 1.194123607785736
$ws.Range("E7").Value = 1.092759629463742
$ws.Range("F7").Value = 1.108116937751023
$ws.Range("G7").Value = 34

$ws.Range("B8").Value = 0.09097433578834678
$ws.Range("C8").Value = 0.7080043655050214
$ws.Range("D8").Value = 1.193680625051876
$ws.Range("E8").Value = 1.092556920737714
$ws.Range("F8").Value = 1.105643785089591
$ws.Range("G8").Value = 33

$ws.Range("B9").Value = 0.03775263814022831
$ws.Range("C9").Value = 0.6203147181500017
$ws.Range("D9").Value = 0.600319873274656
$ws.Range("E9").Value = 0.7748031190403508
$ws.Range("F9").Value = 0.7992627346669655
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = 0.02639662487600095
$ws.Range("C10").Value = 0.6628602414482503
$ws.Range("D10").Value = 0.693693318076211
$ws.Range("E10").Value = 0.8328825355812253
$ws.Range("F10").Value = 0.8774942458002957
$ws.Range("G10").Value = 10

$ws.Range("B11").Value = 0.2236099196487949
$ws.Range("C11").Value = 0.5232375926460167
$ws.Range("D11").Value = 0.3429526541974182
$ws.Range("E11").Value = 0.5856215964233373
$ws.Range("F11").Value = 0.6051355819484564
$ws.Range("G11").Value = 5

